$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = '29.906.47'
$ws.Range("E2").Value = '  +0.12%  '

$ws.Range("D3").Value = '1.905.22'
$ws.Range("E3").Value = '  +0.47%  '

$ws.Range("E4").Value = '  +0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7981'
$ws.Range("E5").Value = '  +5.42%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '240.89'
$ws.Range("E6").Value = '  +0.27%  '

$ws.Range("E7").Value = '  +0.13%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3126'

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '26.15'
$ws.Range("E9").Value = '  +3.19%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06873'
$ws.Range("E10").Value = '  +0.45%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07985'
$ws.Range("E11").Value = '  +0.00%  '

$ws.Range("D12").Value = '1.915.13'
$ws.Range("E12").Value = '  +0.75%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.7337'
$ws.Range("E13").Value = '  -2.10%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.170'
$ws.Range("E14").Value = '  -0.87%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '92.71'
$ws.Range("E15").Value = '  +1.80%  '

$ws.Range("D16").Value = '29.934.47'
$ws.Range("E16").Value = '  +0.16%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.90'
$ws.Range("E17").Value = '  -0.13%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.847'
$ws.Range("E18").Value = '  -1.92%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '244.80'
$ws.Range("E19").Value = '  +2.19%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007697'
$ws.Range("E20").Value = '  -0.26%  '

$ws.Range("D22").Value = '2.155.88'
$ws.Range("E22").Value = '  +0.17%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.003'
$ws.Range("E23").Value = '  +0.14%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.915'
$ws.Range("E24").Value = '  -1.06%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '168.10'
$ws.Range("E25").Value = '  +1.56%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.178'
$ws.Range("E26").Value = '  -0.67%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1414'
$ws.Range("E27").Value = '  +9.07%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.85'
$ws.Range("E28").Value = '  +0.59%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.017'
$ws.Range("E29").Value = '  -0.77%  '

$ws.Range("E30").Value = '  +0.29%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.512'
$ws.Range("E31").Value = '  -0.38%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.284'
$ws.Range("E32").Value = '  -0.15%  '

$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05520'
$ws.Range("E33").Value = '  +2.96%  '

$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.057'
$ws.Range("E34").Value = '  +0.69%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.255'
$ws.Range("E35").Value = '  +0.37%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7283'
$ws.Range("E36").Value = '  +0.07%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.726'
$ws.Range("E37").Value = '  +0.14%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01922'
$ws.Range("E38").Value = '  -0.09%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.790'
$ws.Range("E39").Value = '  +0.83%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.132'
$ws.Range("E40").Value = '  -0.99%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4391'
$ws.Range("E41").Value = '  -0.40%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '72.00'
$ws.Range("E42").Value = '  -0.47%  '

$ws.Range("E43").Value = '  +0.12%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8349'
$ws.Range("E44").Value = '  +0.65%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.865'
$ws.Range("E45").Value = '  -2.77%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '100.49'
$ws.Range("E46").Value = '  -0.49%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.530'
$ws.Range("E47").Value = '  -0.62%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.730'
$ws.Range("E48").Value = '  -0.81%  '

$ws.Range("B49").Value = 'Maker'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '978.11'
$ws.Range("E49").Value = '  +6.86%  '

$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").Value = '2.063.22'
$ws.Range("E50").Value = '  +0.60%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '36.14'
$ws.Range("E51").Value = '  -0.07%  '
